$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find a paragraph index (1-based) whose visible text (paragraph mark
# trimmed) equals $text, searching starting at paragraph $startAt (defaults
# to 1). Returns -1 if not found.
# ---------------------------------------------------------------------------
function Get-ParaIndexByText($doc, $text, $startAt) {
    if (-not $startAt) { $startAt = 1 }
    for ($i = $startAt; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs($i).Range.Text.TrimEnd([char]13)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# Change 1: insert a new bold "Update 4" paragraph right after "Update 3"
# and before "Apps That Are Needed (Gap-Filling Ideas)".
# ---------------------------------------------------------------------------
$idxUpdate3 = Get-ParaIndexByText $d "Update 3"
$p3 = $d.Paragraphs($idxUpdate3)
$p3.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($idxUpdate3 + 1)
$newPara.Range.Text = "Update 4"

# ---------------------------------------------------------------------------
# Changes 2-4: a handful of sentences were previously split across multiple
# runs (separated by <w:proofErr/> grammar-check markers). Re-typing the
# full sentence over the original range merges it back into a single run
# and drops the now orphaned proofErr markers.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    ": Consolidate administrative, training, and operational data for easy access.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ": Consolidate administrative, training, and operational data for easy access.", 2) | Out-Null

$d.Content.Find.Execute(
    "Task tracking categorized by soldier, equipment, or mission.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Task tracking categorized by soldier, equipment, or mission.", 2) | Out-Null

$d.Content.Find.Execute(
    ": Simplify tracking vehicle, weapon, and supply readiness at the platoon level.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ": Simplify tracking vehicle, weapon, and supply readiness at the platoon level.", 2) | Out-Null

# ---------------------------------------------------------------------------
# Changes 5-7: the document's pagination shifted (because of the newly
# added "Update 4" paragraph), so the <w:lastRenderedPageBreak/> markers
# need to move one bullet earlier in three places. These markers are not
# exposed directly on the Word object model, so each affected paragraph is
# rebuilt via Range.InsertXML with the marker added/removed in place.
# ---------------------------------------------------------------------------

function Add-LastRenderedPageBreak($doc, $paraText, $numPr) {
    $idx = Get-ParaIndexByText $doc $paraText
    $p = $doc.Paragraphs($idx)
    $r = $p.Range
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr>' + $numPr + '<w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:lastRenderedPageBreak/><w:t>' + $paraText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml)
}

function Remove-LastRenderedPageBreak($doc, $paraText, $numPr) {
    $idx = Get-ParaIndexByText $doc $paraText
    $p = $doc.Paragraphs($idx)
    $r = $p.Range
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr>' + $numPr + '<w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>' + $paraText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml)
}

function Add-LastRenderedPageBreakToFeatures($doc, $followingText) {
    $idxFollow = Get-ParaIndexByText $doc $followingText
    $idx = $idxFollow - 1
    $p = $doc.Paragraphs($idx)
    $r = $p.Range
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:lastRenderedPageBreak/><w:t>Features</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml)
}

# 5) Section "3. Maintenance and Supply Tracker": the break moves from
#    "Alerts for missing or overdue inspections." to the previous bullet
#    "Digital checklist for scheduled maintenance (linked to GCSS-Army)."
$numPrLvl1Num3 = '<w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr>'
Add-LastRenderedPageBreak $d "Digital checklist for scheduled maintenance (linked to GCSS-Army)." $numPrLvl1Num3
Remove-LastRenderedPageBreak $d "Alerts for missing or overdue inspections." $numPrLvl1Num3

# 6) Section "6. Time Management and Scheduling App": the break moves from
#    "Calendar for meetings, training events, and personal reminders." to
#    the "Features" heading paragraph right before it.
$numPrLvl1Num6 = '<w:numPr><w:ilvl w:val="1"/><w:numId w:val="6"/></w:numPr>'
Add-LastRenderedPageBreakToFeatures $d "Calendar for meetings, training events, and personal reminders."
Remove-LastRenderedPageBreak $d "Calendar for meetings, training events, and personal reminders." $numPrLvl1Num6

# 7) Section "9. Team Climate and Morale Monitor": the break moves from
#    "Anonymous morale surveys soldiers can complete on mobile devices." to
#    the "Features" heading paragraph right before it.
$numPrLvl1Num9 = '<w:numPr><w:ilvl w:val="1"/><w:numId w:val="9"/></w:numPr>'
Add-LastRenderedPageBreakToFeaturesNum9 $d
Remove-LastRenderedPageBreak $d "Anonymous morale surveys soldiers can complete on mobile devices." $numPrLvl1Num9
